# Atualizacao de bases das ligas, do dia: 03-04-2024 as 22:09
# Corrects mismatched match rows (id/teams/odds got paired with the wrong
# fixture) by fixing the data for rows 112, 113, 118, 119, 137, 140 and
# refreshing a handful of odds-only cells in rows 147, 149 and 150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 112
$ws.Range("B112").Value = 6779678
$ws.Range("F112").Value = "Hvidovre IF"
$ws.Range("G112").Value = "Randers FC"
$ws.Range("H112").Value = 1
$ws.Range("I112").Value = 3
$ws.Range("J112").Value = "A"
$ws.Range("K112").Value = 3.3
$ws.Range("L112").Value = 3.6
$ws.Range("M112").Value = 2.05
$ws.Range("N112").Value = 4.5
$ws.Range("O112").Value = 3.8
$ws.Range("P112").Value = 1.8
$ws.Range("Q112").Value = 0.75
$ws.Range("T112").Value = 2.5
$ws.Range("U112").Value = 1.9
$ws.Range("V112").Value = 1.95
$ws.Range("W112").Value = -1
$ws.Range("Y112").Value = 0.8
$ws.Range("Z112").Value = -1
$ws.Range("AA112").Value = 1
$ws.Range("AB112").Value = 0.8999999999999999

# Row 113
$ws.Range("B113").Value = 6779681
$ws.Range("F113").Value = "FC Nordsjaelland"
$ws.Range("G113").Value = "Lyngby"
$ws.Range("H113").Value = 3
$ws.Range("I113").Value = 2
$ws.Range("J113").Value = "H"
$ws.Range("K113").Value = 1.533
$ws.Range("L113").Value = 4.333
$ws.Range("M113").Value = 5
$ws.Range("N113").Value = 1.4
$ws.Range("O113").Value = 5.5
$ws.Range("P113").Value = 6.5
$ws.Range("Q113").Value = -1.25
$ws.Range("T113").Value = 3
$ws.Range("U113").Value = 1.95
$ws.Range("V113").Value = 1.9
$ws.Range("W113").Value = 0.3999999999999999
$ws.Range("Y113").Value = -1
$ws.Range("Z113").Value = -0.5
$ws.Range("AA113").Value = 0.5
$ws.Range("AB113").Value = 0.95

# Row 118
$ws.Range("B118").Value = 6780974
$ws.Range("F118").Value = "Hvidovre IF"
$ws.Range("G118").Value = "Viborg"
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = "D"
$ws.Range("K118").Value = 4
$ws.Range("L118").Value = 3.7
$ws.Range("M118").Value = 1.833
$ws.Range("N118").Value = 5.25
$ws.Range("O118").Value = 4.2
$ws.Range("P118").Value = 1.615
$ws.Range("Q118").Value = 1
$ws.Range("R118").Value = 1.8
$ws.Range("S118").Value = 2.05
$ws.Range("T118").Value = 2.75
$ws.Range("U118").Value = 2
$ws.Range("V118").Value = 1.85
$ws.Range("W118").Value = -1
$ws.Range("X118").Value = 3.2
$ws.Range("Z118").Value = 0.8
$ws.Range("AB118").Value = 1
$ws.Range("AC118").Value = -1

# Row 119
$ws.Range("B119").Value = 6779686
$ws.Range("F119").Value = "Vejle"
$ws.Range("G119").Value = "Silkeborg IF"
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = "H"
$ws.Range("K119").Value = 3
$ws.Range("L119").Value = 3.4
$ws.Range("M119").Value = 2.3
$ws.Range("N119").Value = 3
$ws.Range("O119").Value = 3.4
$ws.Range("P119").Value = 2.3
$ws.Range("Q119").Value = 0.25
$ws.Range("R119").Value = 1.85
$ws.Range("S119").Value = 2
$ws.Range("T119").Value = 2.25
$ws.Range("U119").Value = 1.825
$ws.Range("V119").Value = 2.025
$ws.Range("W119").Value = 2
$ws.Range("X119").Value = -1
$ws.Range("Z119").Value = 0.8500000000000001
$ws.Range("AB119").Value = -0.5
$ws.Range("AC119").Value = 0.5125

# Row 137
$ws.Range("B137").Value = 6779697
$ws.Range("F137").Value = "AGF Aarhus"
$ws.Range("G137").Value = "Hvidovre IF"
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 1.363
$ws.Range("L137").Value = 4.75
$ws.Range("M137").Value = 7.5
$ws.Range("N137").Value = 1.533
$ws.Range("O137").Value = 4.2
$ws.Range("P137").Value = 6
$ws.Range("Q137").Value = -1
$ws.Range("R137").Value = 1.875
$ws.Range("S137").Value = 1.975
$ws.Range("U137").Value = 1.925
$ws.Range("V137").Value = 1.925
$ws.Range("W137").Value = 0.5329999999999999
$ws.Range("Z137").Value = 0
$ws.Range("AA137").Value = -0
$ws.Range("AB137").Value = -1
$ws.Range("AC137").Value = 0.925

# Row 140
$ws.Range("B140").Value = 6779698
$ws.Range("F140").Value = "Brondby"
$ws.Range("G140").Value = "Silkeborg IF"
$ws.Range("H140").Value = 4
$ws.Range("I140").Value = 1
$ws.Range("K140").Value = 1.615
$ws.Range("L140").Value = 3.8
$ws.Range("M140").Value = 5.25
$ws.Range("N140").Value = 1.615
$ws.Range("O140").Value = 3.8
$ws.Range("P140").Value = 5.75
$ws.Range("Q140").Value = -0.75
$ws.Range("R140").Value = 1.85
$ws.Range("S140").Value = 2.05
$ws.Range("U140").Value = 1.85
$ws.Range("V140").Value = 2
$ws.Range("W140").Value = 0.615
$ws.Range("Z140").Value = 0.8500000000000001
$ws.Range("AA140").Value = -1
$ws.Range("AB140").Value = 0.8500000000000001
$ws.Range("AC140").Value = -1

# Row 147
$ws.Range("N147").Value = 3
$ws.Range("P147").Value = 2.375
$ws.Range("Q147").Value = 0.25
$ws.Range("R147").Value = 1.84
$ws.Range("S147").Value = 2.06
$ws.Range("U147").Value = 1.85
$ws.Range("V147").Value = 2

# Row 149
$ws.Range("R149").Value = 1.87
$ws.Range("S149").Value = 2.03

# Row 150
$ws.Range("O150").Value = 3.6
$ws.Range("R150").Value = 1.84
$ws.Range("S150").Value = 2.06
$ws.Range("U150").Value = 1.9
$ws.Range("V150").Value = 1.95
